# feat: add 2022-Q3 data
#
# The workbook originally has 2 sheets: "总计" (summary) and "2021-Q1" (quarterly
# holding detail). This change adds a new "2022-Q3" quarterly sheet (inserted
# right after "总计", before "2021-Q1"), fills it with the latest holding data,
# and records the new quarter as a new row on the "总计" sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "总计" summary sheet
$ws2 = $wb.Worksheets.Item(2)   # currently "2021-Q1" quarterly detail sheet

# 1) Duplicate the existing "2021-Q1" sheet so its data/formatting survive
#    untouched as its own tab; the duplicate is placed right after $ws2 and
#    becomes the active tab (matches original sheet behaviour).
$ws2.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)

# 2) Turn the original $ws2 into the new "2022-Q3" sheet with updated
#    figures, freeing up the "2021-Q1" name for the duplicate.
$ws2.Name = "2022-Q3"
$ws3.Name = "2021-Q1"

$ws2.Range("D1").Value = "基金规模"

# The refreshed figures are stored as text (matching the source data feed),
# not numbers - force text entry, then drop back to the default "Normal"
# style so these cells keep the workbook's plain (unstyled) look.
$ws2.Range("D2:G2").NumberFormat = "@"
$ws2.Range("D2").Value = "0.41"
$ws2.Range("E2").Value = "91.47"
$ws2.Range("F2").Value = "1.82"
$ws2.Range("G2").Value = "0.0075"
$ws2.Range("D2:G2").Style = "Normal"

# Match the header/first-data-row styling used on the "总计" sheet (style
# index 2), copying it onto the refreshed "2022-Q3" sheet.
$ws1.Range("B1:D1").Copy()
$ws2.Range("B1:H1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

# Match the "总计" sheet's page margins on the new "2022-Q3" sheet.
$ws2.PageSetup.LeftMargin = $ws1.PageSetup.LeftMargin
$ws2.PageSetup.RightMargin = $ws1.PageSetup.RightMargin
$ws2.PageSetup.TopMargin = $ws1.PageSetup.TopMargin
$ws2.PageSetup.BottomMargin = $ws1.PageSetup.BottomMargin
$ws2.PageSetup.HeaderMargin = $ws1.PageSetup.HeaderMargin
$ws2.PageSetup.FooterMargin = $ws1.PageSetup.FooterMargin

# 3) Record the new quarter on the "总计" summary sheet: rename the existing
#    row to "2022-Q3" and append a fresh "2021-Q1" row below it.
$ws1.Range("B2").Value = "2022-Q3"

$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "2021-Q1"
$ws1.Range("C3").Value = 1
$ws1.Range("D3").Value = 0.01

$ws1.Range("A2").Copy()
$ws1.Range("A3").PasteSpecial(-4122)
